$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Trends Status")
$ws.Range("B2").Value = 98
$ws.Range("D2").Value = 29
$ws.Range("E2").Value = 17.8
$ws.Range("B3").Value = 106
$ws.Range("D3").Value = 31.4
$ws.Range("E3").Value = 21.7
$ws.Range("B4").Value = 98
$ws.Range("C4").Value = 189
$ws.Range("D4").Value = 29
$ws.Range("E4").Value = 52.6
$ws.Range("B5").Value = 19
$ws.Range("D5").Value = 5.6
$ws.Range("E5").Value = 4.7
$ws.Range("B6").Value = 17
$ws.Range("C6").Value = 11
$ws.Range("D6").Value = 5
$ws.Range("E6").Value = 3.1
$ws.Range("B7").Value = 185
$ws.Range("C7").Value = 284

$ws = $wb.Worksheets.Item("SoIB vs IUCN (SoIB %)")
$ws.Range("B2").Value = 7.9
$ws.Range("C2").Value = 8.4
$ws.Range("D2").Value = 23.6
$ws.Range("E2").Value = 9.6
$ws.Range("F2").Value = 50.6
$ws.Range("E3").Value = 12.1
$ws.Range("F3").Value = 83.3
$ws.Range("D4").Value = 0.5
$ws.Range("E4").Value = 2.5
$ws.Range("F4").Value = 95.7

$ws = $wb.Worksheets.Item("Reason for uplisting")
$ws.Range("B2").Value = 22
$ws.Range("C2").Value = 25
$ws.Range("C4").Value = 36.4
$ws.Range("C5").Value = 2.3
$ws.Range("B6").Value = 20
$ws.Range("C6").Value = 22.7
$ws.Range("C7").Value = 1.1
$ws.Range("C8").Value = 1.1
$ws.Range("C9").Value = 11.4

$ws = $wb.Worksheets.Item("Reason for downlisting")
$ws.Range("B2").Value = 4
$ws.Range("C2").Value = 16
$ws.Range("B3").Value = 5
$ws.Range("C3").Value = 20

$ws = $wb.Worksheets.Item("Priority Status")
$ws.Range("B2").Value = 178
$ws.Range("B3").Value = 323
$ws.Range("B4").Value = 441

$ws = $wb.Worksheets.Item("Species qualification")
$ws.Range("C3").Value = 338
$ws.Range("C4").Value = 359

$ws = $wb.Worksheets.Item("High Priority break-up")
$ws.Range("B2").Value = 77
$ws.Range("C2").Value = 43.3
$ws.Range("D2").Value = 71
$ws.Range("E2").Value = 68.3
$ws.Range("B3").Value = 17
$ws.Range("C3").Value = 9.6
$ws.Range("D3").Value = 7
$ws.Range("E3").Value = 6.7
$ws.Range("C4").Value = 25.3
$ws.Range("E4").Value = 21.2
$ws.Range("C5").Value = 21.9
$ws.Range("E5").Value = 3.8

$ws = $wb.Worksheets.Item("SoIB 2020 vs 2023")
$ws.Range("C3").Value = 3
$ws.Range("D3").Value = 3
$ws.Range("C4").Value = 22
$ws.Range("D4").Value = 21.8
$ws.Range("C7").Value = 333
$ws.Range("D7").Value = 75.5
$ws.Range("C8").Value = 72
$ws.Range("D8").Value = 16.3
$ws.Range("C10").Value = 56
$ws.Range("D10").Value = 17.6
$ws.Range("C11").Value = 81
$ws.Range("D11").Value = 25.4
$ws.Range("C12").Value = 166
$ws.Range("D12").Value = 52

$ws = $wb.Worksheets.Item("SoIB 2023 vs 2020")
$ws.Range("D2").Value = 41.6
$ws.Range("D3").Value = 18
$ws.Range("C4").Value = 56
$ws.Range("D4").Value = 31.5
$ws.Range("D5").Value = 9
$ws.Range("C6").Value = 3
$ws.Range("D6").Value = 0.7
$ws.Range("C7").Value = 333
$ws.Range("D7").Value = 75.5
$ws.Range("C8").Value = 81
$ws.Range("D8").Value = 18.4
$ws.Range("D9").Value = 5.4
$ws.Range("C10").Value = 22
$ws.Range("D10").Value = 6.8
$ws.Range("C11").Value = 72
$ws.Range("D11").Value = 22.3
$ws.Range("C12").Value = 166
$ws.Range("D12").Value = 51.4
$ws.Range("D13").Value = 19.5

$ws = $wb.Worksheets.Item("SoIB vs IUCN (no.)")
$ws.Range("B6").Value = 90
$ws.Range("C6").Value = 269
$ws.Range("D6").Value = 422
$ws.Range("B8").Value = 178
$ws.Range("C8").Value = 323
$ws.Range("D8").Value = 441

$ws = $wb.Worksheets.Item("SoIB vs IUCN (IUCN %)")
$ws.Range("B6").Value = 11.5
$ws.Range("C6").Value = 34.4
$ws.Range("D6").Value = 54

Write-Output "edit complete"
